$wb = $excel.ActiveWorkbook

# ---------- Sheet: SCHEME_MEASURES ----------
$wsScheme = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsScheme.Range("A2").Value = "MQME001"
$wsScheme.Range("A3").Value = "MQME002"
$wsScheme.Range("A4").Value = "MQME003"
$wsScheme.Range("A5").Value = "MQME004"
$wsScheme.Range("A6").Value = "MQME005"

Write-Host "done scheme_measures"

# ---------- Sheet: METADATA_ISSUES ----------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
for ($r = 2; $r -le 7; $r++) {
    $wsIssues.Range("A$r").Value = "MQME012"
}
$wsIssues.Range("A8").Value = "MQME015"
for ($r = 9; $r -le 79; $r++) {
    $wsIssues.Range("A$r").Value = "MQME008"
}

Write-Host "done metadata_issues"

# ---------- Sheet: METADATA_MEASURES ----------
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
# remove the old "Total number of columns" row; rows 3 & 4 shift up to become 2 & 3
$wsMeasures.Rows.Item(2).Delete()
$wsMeasures.Range("A2").Value = "MQME006"
$wsMeasures.Range("A3").Value = "MQME007"

Write-Host "done metadata_measures"

# ---------- Sheet: METADATA_METRICS ----------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")
$metricsData = @(
    ,@("MQID001", "Table names in singular", "94.00%")
    ,@("MQID002", "Table with recommended name length", "100.00%")
    ,@("MQID003", "Columns with correct prefixes", "100.00%")
    ,@("MQID004", "Columns with recommended name size", "99.88%")
    ,@("MQID005", "Columns with comments", "91.22%")
    ,@("MQID006", "Table with standard PK prefixes", "100.00%")
    ,@("MQID007", "Table with standard FK prefixes", "100.00%")
    ,@("MQID008", "Table with standard UK prefixes", "100.00%")
    ,@("MQID009", "NUMBER columns with valid scale", "100.00%")
    ,@("MQID010", "Columns with valid num_distinct", "100.00%")
    ,@("MQID011", "Columns with valid num_nulls", "100.00%")
)

# Force the Value column to stay textual (avoid Excel auto-converting "94.00%" to a percentage number)
$wsMetrics.Range("C2:C12").NumberFormat = "@"

$r = 2
foreach ($row in $metricsData) {
    $wsMetrics.Range("A$r").Value = $row[0]
    $wsMetrics.Range("B$r").Value = $row[1]
    $wsMetrics.Range("C$r").Value = $row[2]
    $r = $r + 1
}

# Drop the temporary text formatting so the cells end up unstyled, like the source data
$wsMetrics.Range("C2:C12").ClearFormats()

Write-Host "done metadata_metrics"
